$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 13.038
$ws.Range("D3").Value = -7.327000000000001
$ws.Range("B4").Value = 5.44
$ws.Range("C4").Value = -13.36
$ws.Range("D4").Value = -8.007999999999999
$ws.Range("C5").Value = -12.78
$ws.Range("E5").Value = 13.113
$ws.Range("B6").Value = 6.164
$ws.Range("C6").Value = -12.089
$ws.Range("B7").Value = 6.234999999999999
$ws.Range("B8").Value = 5.727
$ws.Range("C8").Value = -12.684
$ws.Range("D9").Value = -8.015000000000001
$ws.Range("D11").Value = -7.637
$ws.Range("D14").Value = -8.028
$ws.Range("B16").Value = 5.954
$ws.Range("C16").Value = -12.59
$ws.Range("D18").Value = -7.702
$ws.Range("B20").Value = 6.064
$ws.Range("E20").Value = 13.229
$ws.Range("B21").Value = 5.809
$ws.Range("C22").Value = -12.521
$ws.Range("D25").Value = -7.637
